$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("model_12_1_4", 0.05706809978282501, -0.2385899554003357, -0.5316265055308003, -0.2576982911964609, 1.043547749519348, 1.201628088951111, 2.787129640579224, 1.947746276855469),
    @("model_12_1_0", 0.06125827019749119, -0.711154987868273, 0.3521551909830088, 0.07426661875105511, 1.038910508155823, 1.660090804100037, 1.178895354270935, 1.433645963668823),
    @("model_12_1_3", 0.06217939703634112, -0.1114125475159635, -0.5464531250386404, -0.2237182695847515, 1.037891149520874, 1.078245759010315, 2.814109802246094, 1.895123243331909),
    @("model_12_1_12", 0.1247654212951732, -0.6629452876825075, -1.073801292395252, -0.6982333851684841, 0.9686266779899597, 1.613319754600525, 3.773735284805298, 2.629985332489014),
    @("model_12_1_1", 0.1277428917442136, -0.4050006180737145, 0.240882348078524, 0.1142738973299683, 0.9653316140174866, 1.363072752952576, 1.381380677223206, 1.3716881275177),
    @("model_12_1_5", 0.135474801854099, -0.1519996695147743, -0.4300895089357668, -0.1728353873443647, 0.9567746520042419, 1.117621779441833, 2.602360963821411, 1.816322922706604),
    @("model_12_1_2", 0.1368072471376177, -0.2872953933976587, 0.03021620822131799, 0.03682190616795744, 0.955299973487854, 1.248879909515381, 1.764733791351318, 1.491634964942932),
    @("model_12_1_6", 0.1663405178782422, -0.1068072814689311, -0.3665352852224253, -0.1227046413699104, 0.9226154685020447, 1.073777914047241, 2.486709833145142, 1.73868727684021),
    @("model_12_1_13", 0.1901471681208923, -0.6999222549835651, -1.088921379124801, -0.7188582026061312, 0.8962684273719788, 1.649193286895752, 3.801249742507935, 2.66192626953125),
    @("model_12_1_14", 0.1915101578093853, -0.6933119150229725, -1.108573313161707, -0.7275320788615722, 0.8947600126266479, 1.642780184745789, 3.837010383605957, 2.67535924911499),
    @("model_12_1_15", 0.1939929379526619, -0.6938248827483873, -1.10513797094858, -0.7258028872451787, 0.8920122981071472, 1.643277883529663, 3.830759048461914, 2.672681331634521),
    @("model_12_1_16", 0.2022584373102928, -0.6403938210858167, -1.103493234706655, -0.707172670267719, 0.8828648924827576, 1.59144127368927, 3.827766180038452, 2.643829584121704),
    @("model_12_1_17", 0.2052589386366718, -0.6304939651881261, -1.096736829509677, -0.700153183072149, 0.879544198513031, 1.581836938858032, 3.815471410751343, 2.632958650588989),
    @("model_12_1_22", 0.211204902002747, -0.5949818488752532, -1.093215378226855, -0.6864287113711791, 0.8729638457298279, 1.54738450050354, 3.809063673019409, 2.611703872680664),
    @("model_12_1_23", 0.2149044092285147, -0.5749722584246013, -1.089496608246132, -0.6777362034908305, 0.8688694834709167, 1.527971982955933, 3.80229640007019, 2.598242044448853),
    @("model_12_1_24", 0.2150704207508614, -0.5709183555941366, -1.091293352203961, -0.6773854356890854, 0.8686857223510742, 1.524039149284363, 3.805566072463989, 2.597698926925659),
    @("model_12_1_19", 0.2190467877082372, -0.5468734627972971, -1.088572402842591, -0.667905696708164, 0.8642851114273071, 1.500711798667908, 3.800614833831787, 2.58301830291748),
    @("model_12_1_21", 0.2208006738004115, -0.5369130898610892, -1.087171061605665, -0.6638274432567095, 0.862343966960907, 1.491048574447632, 3.798064231872559, 2.576702356338501),
    @("model_12_1_18", 0.221299118031314, -0.5368051273747769, -1.084920192491017, -0.6625475086439476, 0.8617924451828003, 1.490943908691406, 3.793968677520752, 2.574720144271851),
    @("model_12_1_20", 0.22267529413147, -0.5195304006809676, -1.089549712934801, -0.6593782750946686, 0.8602692484855652, 1.474184632301331, 3.802392959594727, 2.569812059402466),
    @("model_12_1_9", 0.2306393311558016, -0.2889210120454686, -0.5673730446289806, -0.2941566868331036, 0.8514556288719177, 1.250457167625427, 2.852178335189819, 2.004208087921143),
    @("model_12_1_11", 0.2353821424790086, -0.6275130350024014, -0.5006836481565826, -0.3695745628052454, 0.8462066650390625, 1.57894492149353, 2.730822086334229, 2.121004819869995),
    @("model_12_1_10", 0.2458224838007791, -0.4173551063424195, -0.5089721198807835, -0.3044587401506924, 0.8346522450447083, 1.375058531761169, 2.745905160903931, 2.020162582397461),
    @("model_12_1_7", 0.2609609680673114, -0.2138954441300116, -0.2239110455247231, -0.07935577793282178, 0.8178984522819519, 1.177670478820801, 2.227174043655396, 1.671554565429688),
    @("model_12_1_8", 0.2614415255701, -0.2191571376706682, -0.2289419862623232, -0.08388255121584876, 0.8173666000366211, 1.182775020599365, 2.236328840255737, 1.67856502532959)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $rowData[$j]
    }
}
